$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New player roster data replacing the old rows 2-19 (A:C)
$data = @(
    @("Jamal Murray", "PG,SG", "Denver Nuggets"),
    @("Kyshawn George", "SG,SF", "Washington Wizards"),
    @("Cason Wallace", "PG,SG", "Oklahoma City Thunder"),
    @("Shai Gilgeous-Alexander", "PG,SG", "Oklahoma City Thunder"),
    @("Bilal Coulibaly", "SG,SF", "Washington Wizards"),
    @("Kyle Kuzma", "SF,PF", "Milwaukee Bucks"),
    @("Zach LaVine", "SG,SF", "Sacramento Kings"),
    @("Naji Marshall", "SG,SF", "Dallas Mavericks"),
    @("Bam Adebayo", "PF,C", "Miami Heat"),
    @("Kris Dunn", "PG,SG", "LA Clippers"),
    @("Khris Middleton", "SF", "Washington Wizards"),
    @("Kyrie Irving", "PG,SG", "Dallas Mavericks"),
    @("Tobias Harris", "SF,PF", "Detroit Pistons"),
    @("John Collins", "PF,C", "Utah Jazz"),
    @("CJ McCollum", "PG,SG", "New Orleans Pelicans"),
    @("Jordan Clarkson", "SG,SF", "Utah Jazz"),
    @("Lauri Markkanen", "SF,PF", "Utah Jazz"),
    @("Jordan Poole", "PG,SG", "Washington Wizards")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row = $row + 1
}
